# 3D Printing Instructions - update support angle / raft layer values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: Support Angle (H) 0 -> 15, Raft Layers (I) 1 -> 2
$ws.Range("H3").Value = 15
$ws.Range("I3").Value = 2

# Row 4: Infill (G) 0.2 -> 0.16, Support Angle (H) 25 -> 0, Raft Layers (I) 3 -> 2
$ws.Range("G4").Value = 0.16
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 2

# Rows 5-19, 22, 23, 25, 26: Infill (H) 25 -> 20
$rows = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,22,23,25,26)
foreach ($r in $rows) {
    $ws.Range("H$r").Value = 20
}

# Update view: scroll position and selection
$ws.Range("H28").Select()

$excel.ActiveWindow.ScrollRow = 10
